$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.068.53'
$ws.Range('E2').Value = '  +6.56%  '
$ws.Range('D3').Value = '3.531.39'
$ws.Range('E3').Value = '  +8.96%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '193.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +10.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '560.12'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.55%  '
$ws.Range('D7').Value = '3.527.90'
$ws.Range('E7').Value = '  +8.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.612'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.643'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.95%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.82'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.88%  '
$ws.Range('E12').Value = '  +16.10%  '
$ws.Range('E13').Value = '  +8.23%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.54'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +6.83%  '
$ws.Range('D15').Value = '4.102.90'
$ws.Range('E15').Value = '  +8.77%  '
$ws.Range('D16').Value = '3.539.24'
$ws.Range('E16').Value = '  +8.74%  '
$ws.Range('E17').Value = '  +6.38%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.47'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +7.62%  '
$ws.Range('D19').Value = '67.269.34'
$ws.Range('E19').Value = '  +6.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.98'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.12%  '
$ws.Range('E21').Value = '  +4.82%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '409.89'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +12.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +7.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.93'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +7.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.26'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.26'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +13.31%  '
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '12.12'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.97'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +9.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.66'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +8.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '679.20'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.92%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.79'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.90%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.88'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.00%  '
$ws.Range('E35').Value = '  +7.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '60.53'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '39.26'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.85%  '
$ws.Range('E38').Value = '  +17.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.398'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.36%  '
$ws.Range('E41').Value = '  +13.32%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.39'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +20.49%  '
$ws.Range('E43').Value = '  +19.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.998'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.29%  '
$ws.Range('D45').Value = '3.026.79'
$ws.Range('E45').Value = '  +4.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.67'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +7.66%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +12.18%  '
$ws.Range('E48').Value = '  +8.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.12'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E50').Value = '  +4.29%  '
